$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / name updates (not numeric-looking, safe to assign directly)
$ws.Range('D2').Value = '71.161.68'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.813.56'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('D7').Value = '3.812.91'
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  +6.12%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  -2.06%  '
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = '4.456.52'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').Value = '3.784.77'
$ws.Range('E16').Value = '  -5.58%  '
$ws.Range('D17').Value = '71.105.43'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('E25').Value = '  -4.11%  '
$ws.Range('D26').Value = '3.962.67'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('E30').Value = '  -4.83%  '
$ws.Range('E31').Value = '  -3.71%  '
$ws.Range('E32').Value = '  -0.46%  '
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('E35').Value = '  -5.52%  '
$ws.Range('D36').Value = '3.778.63'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  -1.43%  '
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E51').Value = '  -2.45%  '

# Numeric-looking text updates: force Text format so Excel keeps them as strings
# (matches source workbook where these are stored as text, e.g. "71.00" not 71)
$numericCells = @('D4','D5','D6','D9','D11','D12','D13','D14','D18','D19','D20','D21','D22','D23','D24','D25','D27','D28','D30','D31','D33','D34','D38','D41','D42','D43','D46','D47','D48','D49','D50','D51')
foreach ($addr in $numericCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '705.82'
$ws.Range('D6').Value = '171.00'
$ws.Range('D9').Value = '0.523'
$ws.Range('D11').Value = '7.66'
$ws.Range('D12').Value = '0.459'
$ws.Range('D13').Value = '0.0000252'
$ws.Range('D14').Value = '35.79'
$ws.Range('D18').Value = '17.48'
$ws.Range('D19').Value = '7.12'
$ws.Range('D20').Value = '0.114'
$ws.Range('D21').Value = '500.74'
$ws.Range('D22').Value = '10.68'
$ws.Range('D23').Value = '0.723'
$ws.Range('D24').Value = '84.27'
$ws.Range('D25').Value = '0.0000144'
$ws.Range('D27').Value = '12.06'
$ws.Range('D28').Value = '10.35'
$ws.Range('D30').Value = '2.03'
$ws.Range('D31').Value = '3.04'
$ws.Range('D33').Value = '7.35'
$ws.Range('D34').Value = '29.05'
$ws.Range('D38').Value = '9.07'
$ws.Range('D41').Value = '1.02'
$ws.Range('D42').Value = '5.94'
$ws.Range('D43').Value = '3.27'
$ws.Range('D46').Value = '167.02'
$ws.Range('D47').Value = '0.000316'
$ws.Range('D48').Value = '49.14'
$ws.Range('D49').Value = '421.11'
$ws.Range('D50').Value = '8.61'
$ws.Range('D51').Value = '0.293'

# Restore default (Normal) style so cells match original unstyled formatting
foreach ($addr in $numericCells) { $ws.Range($addr).Style = "Normal" }
